$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Team name updates (H = home_team, I = away_team) ---
$ws.Range("H2").Value = "Denver"
$ws.Range("I2").Value = "Brooklyn"

$ws.Range("H3").Value = "Charlotte"
$ws.Range("I3").Value = "Cleveland"

$ws.Range("H4").Value = "Philadelphia"
$ws.Range("I4").Value = "Washington"

$ws.Range("H5").Value = "NewOrleans"
$ws.Range("I5").Value = "Portland"

$ws.Range("H6").Value = "SanAntonio"
$ws.Range("I6").Value = "OklahomaCity"

$ws.Range("H7").Value = "LALakers"
$ws.Range("I7").Value = "NewYork"

# --- Row 2 numeric updates ---
$ws.Range("D2").Value = 230.5
$ws.Range("F2").Value = 115.3134328358209
$ws.Range("G2").Value = 9
$ws.Range("J2").Value = 0.4810040705563093
$ws.Range("K2").Value = 97.91716417910447
$ws.Range("L2").Value = 117.3813432835821
$ws.Range("M2").Value = 114.905223880597
$ws.Range("N2").Value = 75.66417910447763
$ws.Range("O2").Value = 0.3766865671641791
$ws.Range("P2").Value = 0.605910447761194
$ws.Range("Q2").Value = 0.2581940298507462
$ws.Range("R2").Value = 12.46268656716418
$ws.Range("S2").Value = 11.73880597014925
$ws.Range("T2").Value = 0.2062873134328358
$ws.Range("U2").Value = 1.007104216906733
$ws.Range("V2").Value = 0.9964331101750119
$ws.Range("W2").Value = 11.1887378492184
$ws.Range("X2").Value = 0.6268656716417911
$ws.Range("Y2").Value = 51
$ws.Range("Z2").Value = 75.8
$ws.Range("AA2").Value = 0.675
$ws.Range("AB2").Value = 0.4981695965774935
$ws.Range("AC2").Value = 1.49
$ws.Range("AD2").Value = 0.01458111551985547
$ws.Range("AE2").Value = 0.01814676616915423
$ws.Range("AF2").Value = 0.4831421635670515

# --- Row 3 numeric updates ---
$ws.Range("D3").Value = 220
$ws.Range("F3").Value = 111.9045183290708
$ws.Range("G3").Value = 8
$ws.Range("J3").Value = 0.4565217391304348
$ws.Range("K3").Value = 97.59555626598464
$ws.Range("L3").Value = 113.6208333333334
$ws.Range("M3").Value = 113.6588341858482
$ws.Range("N3").Value = 76.5110400682012
$ws.Range("O3").Value = 0.3654807118499573
$ws.Range("P3").Value = 0.5712084398976982
$ws.Range("Q3").Value = 0.2664077152600172
$ws.Range("R3").Value = 11.73751065643649
$ws.Range("S3").Value = 12.93559249786871
$ws.Range("T3").Value = 0.2075799232736573
$ws.Range("U3").Value = 0.9773320378084782
$ws.Range("V3").Value = 0.9650809637243728
$ws.Range("W3").Value = 11.07521281995384
$ws.Range("X3").Value = 0.4661125319693095
$ws.Range("Y3").Value = 40.5
$ws.Range("Z3").Value = 76.40000000000001
$ws.Range("AA3").Value = 0.4098039215686274
$ws.Range("AB3").Value = 0.4974140250207192
$ws.Range("AC3").Value = 1.285
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0.5140153480374476

# --- Row 4 numeric updates ---
$ws.Range("D4").Value = 227.5
$ws.Range("F4").Value = 114.157394843962
$ws.Range("G4").Value = 7.5
$ws.Range("J4").Value = 0.5307692307692308
$ws.Range("K4").Value = 97.12055630936229
$ws.Range("L4").Value = 116.9424016282225
$ws.Range("M4").Value = 115.0138738127544
$ws.Range("N4").Value = 76.93602442333787
$ws.Range("O4").Value = 0.3770835594753505
$ws.Range("P4").Value = 0.5982395974672094
$ws.Range("Q4").Value = 0.2867357530529172
$ws.Range("R4").Value = 12.11282225237449
$ws.Range("S4").Value = 11.60120985979195
$ws.Range("T4").Value = 0.2217438941655359
$ws.Range("U4").Value = 0.9970078152311093
$ws.Range("V4").Value = 1.065397554070877
$ws.Range("W4").Value = 11.60329838741389
$ws.Range("X4").Value = 0.5646766169154229
$ws.Range("Y4").Value = 43
$ws.Range("Z4").Value = 76.15000000000001
$ws.Range("AA4").Value = 0.6277777777777778
$ws.Range("AB4").Value = 0.4857859310167468
$ws.Range("AC4").Value = -0.3300000000000001
$ws.Range("AD4").Value = 0
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.4865730634360924

# --- Row 5 numeric updates ---
$ws.Range("D5").Value = 230.5
$ws.Range("F5").Value = 114.3513116236997
$ws.Range("G5").Value = 1
$ws.Range("J5").Value = 0.4847301951779563
$ws.Range("K5").Value = 97.90226142017184
$ws.Range("L5").Value = 116.1303256445048
$ws.Range("M5").Value = 116.4183061962913
$ws.Range("N5").Value = 76.42901402080508
$ws.Range("O5").Value = 0.382421867933062
$ws.Range("P5").Value = 0.5901792175486207
$ws.Range("Q5").Value = 0.2938032564450475
$ws.Range("R5").Value = 12.569855269109
$ws.Range("S5").Value = 12.20900045228403
$ws.Range("T5").Value = 0.2244795906829489
$ws.Range("U5").Value = 0.9987014115606959
$ws.Range("V5").Value = 0.9357113316933782
$ws.Range("W5").Value = 11.47363797512983
$ws.Range("X5").Value = 0.473767526006332
$ws.Range("Y5").Value = 42.5
$ws.Range("Z5").Value = 76.25
$ws.Range("AA5").Value = 0.4852941176470588
$ws.Range("AB5").Value = 0.5029450642138835
$ws.Range("AC5").Value = -4.695
$ws.Range("AD5").Value = 0.1412615710230504
$ws.Range("AE5").Value = 0.1307834957158651
$ws.Range("AF5").Value = 0.457843756037993

# --- Row 6 numeric updates ---
$ws.Range("D6").Value = 237.5
$ws.Range("F6").Value = 115.2424242424242
$ws.Range("G6").Value = 3.5
$ws.Range("J6").Value = 0.5875699626865671
$ws.Range("K6").Value = 100.6174242424242
$ws.Range("L6").Value = 113.6886363636364
$ws.Range("M6").Value = 118.0416666666667
$ws.Range("N6").Value = 74.13712121212122
$ws.Range("O6").Value = 0.3527348484848485
$ws.Range("P6").Value = 0.5662196969696969
$ws.Range("Q6").Value = 0.2430909090909091
$ws.Range("R6").Value = 11.9969696969697
$ws.Range("S6").Value = 12.88787878787879
$ws.Range("T6").Value = 0.2034469696969697
$ws.Range("U6").Value = 1.006484054518989
$ws.Range("V6").Value = 1.034023897491248
$ws.Range("W6").Value = 11.59518142333837
$ws.Range("X6").Value = 0.3636363636363636
$ws.Range("Y6").Value = 23
$ws.Range("Z6").Value = 73.80000000000001
$ws.Range("AA6").Value = 0.3851674641148325
$ws.Range("AB6").Value = 0.5000213014461756
$ws.Range("AC6").Value = -3.94
$ws.Range("AD6").Value = 0.02036293476560462
$ws.Range("AE6").Value = 0.02138415404040404
$ws.Range("AF6").Value = 0.4924064783506635

# --- Row 7 numeric updates ---
$ws.Range("D7").Value = 225.5
$ws.Range("F7").Value = 116.0273266022827
$ws.Range("G7").Value = 2.5
$ws.Range("J7").Value = 0.5149253731343284
$ws.Range("K7").Value = 98.68252853380157
$ws.Range("L7").Value = 116.1953577699737
$ws.Range("M7").Value = 115.0210381913959
$ws.Range("N7").Value = 76.51334503950832
$ws.Range("O7").Value = 0.3730020851624233
$ws.Range("P7").Value = 0.5765320456540826
$ws.Range("Q7").Value = 0.2907157594381036
$ws.Range("R7").Value = 11.25847234416155
$ws.Range("S7").Value = 10.60206321334504
$ws.Range("T7").Value = 0.2104190079016681
$ws.Range("U7").Value = 1.013339096963168
$ws.Range("V7").Value = 1.005596910445933
$ws.Range("W7").Value = 10.87727887773243
$ws.Range("X7").Value = 0.5330333625987709
$ws.Range("Y7").Value = 41.5
$ws.Range("Z7").Value = 75.40000000000001
$ws.Range("AA7").Value = 0.5463659147869674
$ws.Range("AB7").Value = 0.5096070111434863
$ws.Range("AC7").Value = -1.545
$ws.Range("AD7").Value = 0.1289987059742818
$ws.Range("AE7").Value = 0.06672505527915974
$ws.Range("AF7").Value = 0.4796797761337289
